# Applies the "mail setup done, deployment left" edit:
#  - strips the surrounding quotes from every "Single Family" Property Type cell
#  - adds the missing trailing period in row 2's note
#  - splits the long, sentence-joined Notes cells across additional columns
#    (one sentence/clause per cell), tidying up wording along the way

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Property Type column: drop the surrounding quotes everywhere ---
$ws.Range("E2").Value = "Single Family"
$ws.Range("E3").Value = "Single Family"
$ws.Range("E4").Value = "Single Family"
$ws.Range("E5").Value = "Single Family"
$ws.Range("E6").Value = "Single Family"
$ws.Range("E7").Value = "Single Family"
$ws.Range("E8").Value = "Single Family"
$ws.Range("E9").Value = "Single Family"
$ws.Range("E10").Value = "Single Family"

# --- Row 2 Notes: add the missing trailing period ---
$ws.Range("I2").Value = """Good shape and great east location."""

# --- Row 3 Notes: split into I3 / J3 ---
$ws.Range("I3").Value = """Large lot"
$ws.Range("J3").Value = " Needs work"""

# --- Row 4 Notes: split into I4 / J4 / K4 / L4 / M4 ---
$ws.Range("I4").Value = """A/C 2021"
$ws.Range("J4").Value = " New garage door"
$ws.Range("K4").Value = " Roof 2004"
$ws.Range("L4").Value = " New exterior paint"
$ws.Range("M4").Value = " Good shape - needs updating"""

# --- Row 5 Notes: split into I5 / J5 / K5 ---
$ws.Range("I5").Value = """Garage converted legally to 4/3"
$ws.Range("J5").Value = " Screened pool"
$ws.Range("K5").Value = " Needs new roof and updates"""

# --- Row 6 Notes: split into I6 / J6 / K6 / L6 / M6 ---
$ws.Range("I6").Value = """Needs interior rehab"
$ws.Range("J6").Value = " Roof 5 years old"
$ws.Range("K6").Value = " Central AC"
$ws.Range("L6").Value = " 55+ community"
$ws.Range("M6").Value = " Buyer to assume tenant leaving in March"""

# --- Row 7 Notes: split into I7 / J7 / K7 / L7 / M7 / N7 / O7 / P7 / Q7 / R7 ---
$ws.Range("I7").Value = """Carport"
$ws.Range("J7").Value = " Huge driveway"
$ws.Range("K7").Value = " Roof 2022"
$ws.Range("L7").Value = " AC 2016"
$ws.Range("M7").Value = " Water heater 2024"
$ws.Range("N7").Value = " Detached workshop"
$ws.Range("O7").Value = " Corner lot"
$ws.Range("P7").Value = " Fenced in"
$ws.Range("Q7").Value = " Low DOM area"
$ws.Range("R7").Value = " Only active in subdivision"""

# --- Row 8 Notes: split into I8 / J8 / K8 ---
$ws.Range("I8").Value = """Garage"
$ws.Range("J8").Value = " 55+ community"
$ws.Range("K8").Value = " Needs some updates"""

# --- Row 9 Notes: split into I9 / J9 / K9 / L9 / M9 / N9 ---
$ws.Range("I9").Value = """55+ community"
$ws.Range("J9").Value = " All new appliances"
$ws.Range("K9").Value = " Roof inspected with 6 years left"
$ws.Range("L9").Value = " New mini-split HVAC system"
$ws.Range("M9").Value = " HOA $276 monthly"
$ws.Range("N9").Value = " 2 year rental restriction"""

# --- Row 10 Notes: split into I10 / J10 / K10 ---
$ws.Range("I10").Value = """Needs full rehab"
$ws.Range("J10").Value = " On city water/sewer"
$ws.Range("K10").Value = " Roof replaced 6 years ago but leaks"""
